$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.2271293375394322
    "C2" = 0.4921135646687697
    "J2" = 0.0220820189274448
    "P2" = 0.1608832807570978
    "S2" = 0.09779179810725552
    "C3" = 0.0189873417721519
    "J3" = 0.06962025316455696
    "P3" = 0.7025316455696202
    "S3" = 0.2088607594936709
    "J4" = 0.03571428571428571
    "P4" = 0.6964285714285714
    "S4" = 0.2678571428571428
    "B6" = 0.05217391304347826
    "D6" = 0.01304347826086956
    "F6" = 0.04782608695652174
    "J6" = 0.2739130434782608
    "O6" = 0.004347826086956522
    "Q6" = 0.1391304347826087
    "R6" = 0.0782608695652174
    "S6" = 0.391304347826087
    "B7" = 0.163716814159292
    "D7" = 0.02212389380530973
    "F7" = 0.05309734513274336
    "J7" = 0.1150442477876106
    "O7" = 0.01327433628318584
    "Q7" = 0.1725663716814159
    "R7" = 0.05752212389380531
    "S7" = 0.4026548672566372
    "B8" = 0.08571428571428572
    "D8" = 0.01428571428571429
    "E8" = 0.002040816326530612
    "F8" = 0.06326530612244897
    "J8" = 0.1163265306122449
    "O8" = 0.01020408163265306
    "Q8" = 0.1693877551020408
    "R8" = 0.1306122448979592
    "S8" = 0.4081632653061225
    "B9" = 0.1333333333333333
    "D9" = 0.03589743589743589
    "F9" = 0.05128205128205128
    "J9" = 0.1538461538461539
    "O9" = 0.005128205128205128
    "Q9" = 0.1487179487179487
    "R9" = 0.1230769230769231
    "S9" = 0.3487179487179487
    "B10" = 0.09767092411720511
    "D10" = 0.0270473328324568
    "E10" = 0.0007513148009015778
    "F10" = 0.06611570247933884
    "J10" = 0.1134485349361382
    "O10" = 0.009015777610818933
    "Q10" = 0.2291510142749812
    "R10" = 0.1051840721262209
    "S10" = 0.3516153268219384
    "G11" = 0.1507936507936508
    "J11" = 0.09788359788359788
    "K11" = 0.2380952380952381
    "L11" = 0.4708994708994709
    "S11" = 0.04232804232804233
    "G12" = 0.664804469273743
    "J12" = 0.2290502793296089
    "K12" = 0.0111731843575419
    "L12" = 0.0223463687150838
    "S12" = 0.07262569832402235
    "G13" = 0.7205882352941176
    "J13" = 0.2352941176470588
    "S13" = 0.04411764705882353
    "G14" = 0.5
    "J14" = 0.5
    "F15" = 0.03553299492385787
    "H15" = 0.1878172588832487
    "I15" = 0.05583756345177665
    "J15" = 0.3807106598984771
    "K15" = 0.06598984771573604
    "M15" = 0.01015228426395939
    "N15" = 0.005076142131979695
    "O15" = 0.04568527918781726
    "S15" = 0.2131979695431472
    "F16" = 0.005050505050505051
    "H16" = 0.1767676767676768
    "I16" = 0.04040404040404041
    "J16" = 0.4242424242424243
    "K16" = 0.1262626262626263
    "M16" = 0.03535353535353535
    "O16" = 0.0505050505050505
    "S16" = 0.1414141414141414
    "F17" = 0.006147540983606557
    "H17" = 0.1905737704918033
    "I17" = 0.1127049180327869
    "J17" = 0.389344262295082
    "K17" = 0.1024590163934426
    "M17" = 0.04098360655737705
    "O17" = 0.06762295081967214
    "S17" = 0.09016393442622951
    "F18" = 0.007782101167315175
    "H18" = 0.198443579766537
    "I18" = 0.0933852140077821
    "J18" = 0.3891050583657588
    "K18" = 0.09727626459143969
    "M18" = 0.02723735408560311
    "N18" = 0.003891050583657588
    "O18" = 0.0622568093385214
    "S18" = 0.1206225680933852
    "F19" = 0.02777777777777778
    "H19" = 0.2137345679012346
    "I19" = 0.07407407407407407
    "J19" = 0.3518518518518519
    "K19" = 0.1319444444444444
    "M19" = 0.02469135802469136
    "O19" = 0.06327160493827161
    "S19" = 0.1126543209876543
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}

Write-Host "Updated $($values.Count) cells"
